# Auto-generated edit script: update crypto price/volume table (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.795.61'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '2.477.13'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('D4').Value = "'" + '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'" + '319.04'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('D6').Value = "'" + '93.30'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').Value = "'" + '0.553'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = "'" + '0.519'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.33%  '
$ws.Range('D10').Value = "'" + '0.0882'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +11.49%  '
$ws.Range('D11').Value = "'" + '33.37'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +2.87%  '
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = '2.858.72'
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = "'" + '15.70'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '2.460.44'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = "'" + '0.800'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.38%  '
$ws.Range('D18').Value = '41.736.92'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '0.0₃0954'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').Value = "'" + '6.48'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.21%  '
$ws.Range('D21').Value = "'" + '71.23'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').Value = "'" + '11.35'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.62%  '
$ws.Range('D23').Value = "'" + '242.70'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').Value = "'" + '2.77'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.84%  '
$ws.Range('E25').Value = '  +2.83%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = "'" + '25.19'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.89%  '
$ws.Range('E28').Value = '  +0.53%  '
$ws.Range('D29').Value = "'" + '9.78'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('D30').Value = "'" + '36.96'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.81%  '
$ws.Range('D31').Value = "'" + '157.86'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('D32').Value = "'" + '5.53'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.01%  '
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('E34').Value = '  +1.22%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = "'" + '17.47'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.13%  '
$ws.Range('E37').Value = '  +5.64%  '
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('E39').Value = '  +1.80%  '
$ws.Range('E40').Value = '  +1.73%  '
$ws.Range('D41').Value = "'" + '4.04'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('D42').Value = "'" + '2.53'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +7.46%  '
$ws.Range('D43').Value = '2.004.04'
$ws.Range('E43').Value = '  +3.23%  '
$ws.Range('D44').Value = "'" + '19.23'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.58%  '
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('D46').Value = "'" + '2.99'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.51%  '
$ws.Range('E47').Value = '  +5.67%  '
$ws.Range('D48').Value = '2.717.44'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('D49').Value = "'" + '76.96'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +8.13%  '
$ws.Range('D50').Value = "'" + '97.99'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('D51').Value = "'" + '67.57'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.12%  '
